$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 235,4
$data[0,0] = 'codeforiati:category-name'
$data[0,1] = 'codeforiati:category-code'
$data[0,2] = 'codeforiati:group-name'
$data[0,3] = 'codeforiati:group-code'
$data[1,0] = 'Educação, nível não especificado'
$data[1,1] = '111'
$data[1,2] = 'Educação'
$data[1,3] = '110'
$data[2,0] = 'Educação, nível não especificado'
$data[2,1] = '111'
$data[2,2] = 'Educação'
$data[2,3] = '110'
$data[3,0] = 'Educação, nível não especificado'
$data[3,1] = '111'
$data[3,2] = 'Educação'
$data[3,3] = '110'
$data[4,0] = 'Educação, nível não especificado'
$data[4,1] = '111'
$data[4,2] = 'Educação'
$data[4,3] = '110'
$data[5,0] = 'Ensino básico'
$data[5,1] = '112'
$data[5,2] = 'Educação'
$data[5,3] = '110'
$data[6,0] = 'Ensino básico'
$data[6,1] = '112'
$data[6,2] = 'Educação'
$data[6,3] = '110'
$data[7,0] = 'Ensino básico'
$data[7,1] = '112'
$data[7,2] = 'Educação'
$data[7,3] = '110'
$data[8,0] = 'Ensino básico'
$data[8,1] = '112'
$data[8,2] = 'Educação'
$data[8,3] = '110'
$data[9,0] = 'Ensino básico'
$data[9,1] = '112'
$data[9,2] = 'Educação'
$data[9,3] = '110'
$data[10,0] = 'Ensino básico'
$data[10,1] = '112'
$data[10,2] = 'Educação'
$data[10,3] = '110'
$data[11,0] = 'Ensino básico'
$data[11,1] = '112'
$data[11,2] = 'Educação'
$data[11,3] = '110'
$data[12,0] = 'Ensino secundário'
$data[12,1] = '113'
$data[12,2] = 'Educação'
$data[12,3] = '110'
$data[13,0] = 'Ensino secundário'
$data[13,1] = '113'
$data[13,2] = 'Educação'
$data[13,3] = '110'
$data[14,0] = 'Ensino pós-secundário'
$data[14,1] = '114'
$data[14,2] = 'Educação'
$data[14,3] = '110'
$data[15,0] = 'Ensino pós-secundário'
$data[15,1] = '114'
$data[15,2] = 'Educação'
$data[15,3] = '110'
$data[16,0] = 'Saúde, geral'
$data[16,1] = '121'
$data[16,2] = 'Saúde'
$data[16,3] = '120'
$data[17,0] = 'Saúde, geral'
$data[17,1] = '121'
$data[17,2] = 'Saúde'
$data[17,3] = '120'
$data[18,0] = 'Saúde, geral'
$data[18,1] = '121'
$data[18,2] = 'Saúde'
$data[18,3] = '120'
$data[19,0] = 'Saúde, geral'
$data[19,1] = '121'
$data[19,2] = 'Saúde'
$data[19,3] = '120'
$data[20,0] = 'Saúde básica'
$data[20,1] = '122'
$data[20,2] = 'Saúde'
$data[20,3] = '120'
$data[21,0] = 'Saúde básica'
$data[21,1] = '122'
$data[21,2] = 'Saúde'
$data[21,3] = '120'
$data[22,0] = 'Saúde básica'
$data[22,1] = '122'
$data[22,2] = 'Saúde'
$data[22,3] = '120'
$data[23,0] = 'Saúde básica'
$data[23,1] = '122'
$data[23,2] = 'Saúde'
$data[23,3] = '120'
$data[24,0] = 'Saúde básica'
$data[24,1] = '122'
$data[24,2] = 'Saúde'
$data[24,3] = '120'
$data[25,0] = 'Saúde básica'
$data[25,1] = '122'
$data[25,2] = 'Saúde'
$data[25,3] = '120'
$data[26,0] = 'Saúde básica'
$data[26,1] = '122'
$data[26,2] = 'Saúde'
$data[26,3] = '120'
$data[27,0] = 'Saúde básica'
$data[27,1] = '122'
$data[27,2] = 'Saúde'
$data[27,3] = '120'
$data[28,0] = 'Saúde básica'
$data[28,1] = '122'
$data[28,2] = 'Saúde'
$data[28,3] = '120'
$data[29,0] = 'Doenças não transmissíveis (DNT)'
$data[29,1] = '123'
$data[29,2] = 'Saúde'
$data[29,3] = '120'
$data[30,0] = 'Doenças não transmissíveis (DNT)'
$data[30,1] = '123'
$data[30,2] = 'Saúde'
$data[30,3] = '120'
$data[31,0] = 'Doenças não transmissíveis (DNT)'
$data[31,1] = '123'
$data[31,2] = 'Saúde'
$data[31,3] = '120'
$data[32,0] = 'Doenças não transmissíveis (DNT)'
$data[32,1] = '123'
$data[32,2] = 'Saúde'
$data[32,3] = '120'
$data[33,0] = 'Doenças não transmissíveis (DNT)'
$data[33,1] = '123'
$data[33,2] = 'Saúde'
$data[33,3] = '120'
$data[34,0] = 'Doenças não transmissíveis (DNT)'
$data[34,1] = '123'
$data[34,2] = 'Saúde'
$data[34,3] = '120'
$data[35,0] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[35,1] = '130'
$data[35,2] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[35,3] = '130'
$data[36,0] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[36,1] = '130'
$data[36,2] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[36,3] = '130'
$data[37,0] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[37,1] = '130'
$data[37,2] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[37,3] = '130'
$data[38,0] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[38,1] = '130'
$data[38,2] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[38,3] = '130'
$data[39,0] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[39,1] = '130'
$data[39,2] = 'Políticas/programas populacionais e saúde reprodutiva'
$data[39,3] = '130'
$data[40,0] = 'Abastecimento de água e saneamento'
$data[40,1] = '140'
$data[40,2] = 'Abastecimento de água e saneamento'
$data[40,3] = '140'
$data[41,0] = 'Abastecimento de água e saneamento'
$data[41,1] = '140'
$data[41,2] = 'Abastecimento de água e saneamento'
$data[41,3] = '140'
$data[42,0] = 'Abastecimento de água e saneamento'
$data[42,1] = '140'
$data[42,2] = 'Abastecimento de água e saneamento'
$data[42,3] = '140'
$data[43,0] = 'Abastecimento de água e saneamento'
$data[43,1] = '140'
$data[43,2] = 'Abastecimento de água e saneamento'
$data[43,3] = '140'
$data[44,0] = 'Abastecimento de água e saneamento'
$data[44,1] = '140'
$data[44,2] = 'Abastecimento de água e saneamento'
$data[44,3] = '140'
$data[45,0] = 'Abastecimento de água e saneamento'
$data[45,1] = '140'
$data[45,2] = 'Abastecimento de água e saneamento'
$data[45,3] = '140'
$data[46,0] = 'Abastecimento de água e saneamento'
$data[46,1] = '140'
$data[46,2] = 'Abastecimento de água e saneamento'
$data[46,3] = '140'
$data[47,0] = 'Abastecimento de água e saneamento'
$data[47,1] = '140'
$data[47,2] = 'Abastecimento de água e saneamento'
$data[47,3] = '140'
$data[48,0] = 'Abastecimento de água e saneamento'
$data[48,1] = '140'
$data[48,2] = 'Abastecimento de água e saneamento'
$data[48,3] = '140'
$data[49,0] = 'Abastecimento de água e saneamento'
$data[49,1] = '140'
$data[49,2] = 'Abastecimento de água e saneamento'
$data[49,3] = '140'
$data[50,0] = 'Abastecimento de água e saneamento'
$data[50,1] = '140'
$data[50,2] = 'Abastecimento de água e saneamento'
$data[50,3] = '140'
$data[51,0] = 'Governo e sociedade civil, geral'
$data[51,1] = '151'
$data[51,2] = 'Governo e sociedade civil'
$data[51,3] = '150'
$data[52,0] = 'Governo e sociedade civil, geral'
$data[52,1] = '151'
$data[52,2] = 'Governo e sociedade civil'
$data[52,3] = '150'
$data[53,0] = 'Governo e sociedade civil, geral'
$data[53,1] = '151'
$data[53,2] = 'Governo e sociedade civil'
$data[53,3] = '150'
$data[54,0] = 'Governo e sociedade civil, geral'
$data[54,1] = '151'
$data[54,2] = 'Governo e sociedade civil'
$data[54,3] = '150'
$data[55,0] = 'Governo e sociedade civil, geral'
$data[55,1] = '151'
$data[55,2] = 'Governo e sociedade civil'
$data[55,3] = '150'
$data[56,0] = 'Governo e sociedade civil, geral'
$data[56,1] = '151'
$data[56,2] = 'Governo e sociedade civil'
$data[56,3] = '150'
$data[57,0] = 'Governo e sociedade civil, geral'
$data[57,1] = '151'
$data[57,2] = 'Governo e sociedade civil'
$data[57,3] = '150'
$data[58,0] = 'Governo e sociedade civil, geral'
$data[58,1] = '151'
$data[58,2] = 'Governo e sociedade civil'
$data[58,3] = '150'
$data[59,0] = 'Governo e sociedade civil, geral'
$data[59,1] = '151'
$data[59,2] = 'Governo e sociedade civil'
$data[59,3] = '150'
$data[60,0] = 'Governo e sociedade civil, geral'
$data[60,1] = '151'
$data[60,2] = 'Governo e sociedade civil'
$data[60,3] = '150'
$data[61,0] = 'Governo e sociedade civil, geral'
$data[61,1] = '151'
$data[61,2] = 'Governo e sociedade civil'
$data[61,3] = '150'
$data[62,0] = 'Governo e sociedade civil, geral'
$data[62,1] = '151'
$data[62,2] = 'Governo e sociedade civil'
$data[62,3] = '150'
$data[63,0] = 'Governo e sociedade civil, geral'
$data[63,1] = '151'
$data[63,2] = 'Governo e sociedade civil'
$data[63,3] = '150'
$data[64,0] = 'Governo e sociedade civil, geral'
$data[64,1] = '151'
$data[64,2] = 'Governo e sociedade civil'
$data[64,3] = '150'
$data[65,0] = 'Governo e sociedade civil, geral'
$data[65,1] = '151'
$data[65,2] = 'Governo e sociedade civil'
$data[65,3] = '150'
$data[66,0] = 'Governo e sociedade civil, geral'
$data[66,1] = '151'
$data[66,2] = 'Governo e sociedade civil'
$data[66,3] = '150'
$data[67,0] = 'Conflitos, paz e segurança'
$data[67,1] = '152'
$data[67,2] = 'Governo e sociedade civil'
$data[67,3] = '150'
$data[68,0] = 'Conflitos, paz e segurança'
$data[68,1] = '152'
$data[68,2] = 'Governo e sociedade civil'
$data[68,3] = '150'
$data[69,0] = 'Conflitos, paz e segurança'
$data[69,1] = '152'
$data[69,2] = 'Governo e sociedade civil'
$data[69,3] = '150'
$data[70,0] = 'Conflitos, paz e segurança'
$data[70,1] = '152'
$data[70,2] = 'Governo e sociedade civil'
$data[70,3] = '150'
$data[71,0] = 'Conflitos, paz e segurança'
$data[71,1] = '152'
$data[71,2] = 'Governo e sociedade civil'
$data[71,3] = '150'
$data[72,0] = 'Conflitos, paz e segurança'
$data[72,1] = '152'
$data[72,2] = 'Governo e sociedade civil'
$data[72,3] = '150'
$data[73,0] = 'Outros serviços e infraestruturas sociais'
$data[73,1] = '160'
$data[73,2] = 'Outros serviços e infraestruturas sociais'
$data[73,3] = '160'
$data[74,0] = 'Outros serviços e infraestruturas sociais'
$data[74,1] = '160'
$data[74,2] = 'Outros serviços e infraestruturas sociais'
$data[74,3] = '160'
$data[75,0] = 'Outros serviços e infraestruturas sociais'
$data[75,1] = '160'
$data[75,2] = 'Outros serviços e infraestruturas sociais'
$data[75,3] = '160'
$data[76,0] = 'Outros serviços e infraestruturas sociais'
$data[76,1] = '160'
$data[76,2] = 'Outros serviços e infraestruturas sociais'
$data[76,3] = '160'
$data[77,0] = 'Outros serviços e infraestruturas sociais'
$data[77,1] = '160'
$data[77,2] = 'Outros serviços e infraestruturas sociais'
$data[77,3] = '160'
$data[78,0] = 'Outros serviços e infraestruturas sociais'
$data[78,1] = '160'
$data[78,2] = 'Outros serviços e infraestruturas sociais'
$data[78,3] = '160'
$data[79,0] = 'Outros serviços e infraestruturas sociais'
$data[79,1] = '160'
$data[79,2] = 'Outros serviços e infraestruturas sociais'
$data[79,3] = '160'
$data[80,0] = 'Outros serviços e infraestruturas sociais'
$data[80,1] = '160'
$data[80,2] = 'Outros serviços e infraestruturas sociais'
$data[80,3] = '160'
$data[81,0] = 'Outros serviços e infraestruturas sociais'
$data[81,1] = '160'
$data[81,2] = 'Outros serviços e infraestruturas sociais'
$data[81,3] = '160'
$data[82,0] = 'Outros serviços e infraestruturas sociais'
$data[82,1] = '160'
$data[82,2] = 'Outros serviços e infraestruturas sociais'
$data[82,3] = '160'
$data[83,0] = 'Outros serviços e infraestruturas sociais'
$data[83,1] = '160'
$data[83,2] = 'Outros serviços e infraestruturas sociais'
$data[83,3] = '160'
$data[84,0] = 'Transportes e armazenamento'
$data[84,1] = '210'
$data[84,2] = 'Transportes e armazenamento'
$data[84,3] = '210'
$data[85,0] = 'Transportes e armazenamento'
$data[85,1] = '210'
$data[85,2] = 'Transportes e armazenamento'
$data[85,3] = '210'
$data[86,0] = 'Transportes e armazenamento'
$data[86,1] = '210'
$data[86,2] = 'Transportes e armazenamento'
$data[86,3] = '210'
$data[87,0] = 'Transportes e armazenamento'
$data[87,1] = '210'
$data[87,2] = 'Transportes e armazenamento'
$data[87,3] = '210'
$data[88,0] = 'Transportes e armazenamento'
$data[88,1] = '210'
$data[88,2] = 'Transportes e armazenamento'
$data[88,3] = '210'
$data[89,0] = 'Transportes e armazenamento'
$data[89,1] = '210'
$data[89,2] = 'Transportes e armazenamento'
$data[89,3] = '210'
$data[90,0] = 'Transportes e armazenamento'
$data[90,1] = '210'
$data[90,2] = 'Transportes e armazenamento'
$data[90,3] = '210'
$data[91,0] = 'Comunicações'
$data[91,1] = '220'
$data[91,2] = 'Comunicações'
$data[91,3] = '220'
$data[92,0] = 'Comunicações'
$data[92,1] = '220'
$data[92,2] = 'Comunicações'
$data[92,3] = '220'
$data[93,0] = 'Comunicações'
$data[93,1] = '220'
$data[93,2] = 'Comunicações'
$data[93,3] = '220'
$data[94,0] = 'Comunicações'
$data[94,1] = '220'
$data[94,2] = 'Comunicações'
$data[94,3] = '220'
$data[95,0] = 'Política energética'
$data[95,1] = '231'
$data[95,2] = 'Energia'
$data[95,3] = '230'
$data[96,0] = 'Política energética'
$data[96,1] = '231'
$data[96,2] = 'Energia'
$data[96,3] = '230'
$data[97,0] = 'Política energética'
$data[97,1] = '231'
$data[97,2] = 'Energia'
$data[97,3] = '230'
$data[98,0] = 'Política energética'
$data[98,1] = '231'
$data[98,2] = 'Energia'
$data[98,3] = '230'
$data[99,0] = 'Geração de energia, fontes renováveis'
$data[99,1] = '232'
$data[99,2] = 'Energia'
$data[99,3] = '230'
$data[100,0] = 'Geração de energia, fontes renováveis'
$data[100,1] = '232'
$data[100,2] = 'Energia'
$data[100,3] = '230'
$data[101,0] = 'Geração de energia, fontes renováveis'
$data[101,1] = '232'
$data[101,2] = 'Energia'
$data[101,3] = '230'
$data[102,0] = 'Geração de energia, fontes renováveis'
$data[102,1] = '232'
$data[102,2] = 'Energia'
$data[102,3] = '230'
$data[103,0] = 'Geração de energia, fontes renováveis'
$data[103,1] = '232'
$data[103,2] = 'Energia'
$data[103,3] = '230'
$data[104,0] = 'Geração de energia, fontes renováveis'
$data[104,1] = '232'
$data[104,2] = 'Energia'
$data[104,3] = '230'
$data[105,0] = 'Geração de energia, fontes renováveis'
$data[105,1] = '232'
$data[105,2] = 'Energia'
$data[105,3] = '230'
$data[106,0] = 'Geração de energia, fontes renováveis'
$data[106,1] = '232'
$data[106,2] = 'Energia'
$data[106,3] = '230'
$data[107,0] = 'Geração de energia, fontes renováveis'
$data[107,1] = '232'
$data[107,2] = 'Energia'
$data[107,3] = '230'
$data[108,0] = 'Geração de energia, fontes não renováveis'
$data[108,1] = '233'
$data[108,2] = 'Energia'
$data[108,3] = '230'
$data[109,0] = 'Geração de energia, fontes não renováveis'
$data[109,1] = '233'
$data[109,2] = 'Energia'
$data[109,3] = '230'
$data[110,0] = 'Geração de energia, fontes não renováveis'
$data[110,1] = '233'
$data[110,2] = 'Energia'
$data[110,3] = '230'
$data[111,0] = 'Geração de energia, fontes não renováveis'
$data[111,1] = '233'
$data[111,2] = 'Energia'
$data[111,3] = '230'
$data[112,0] = 'Geração de energia, fontes não renováveis'
$data[112,1] = '233'
$data[112,2] = 'Energia'
$data[112,3] = '230'
$data[113,0] = 'Geração de energia, fontes não renováveis'
$data[113,1] = '233'
$data[113,2] = 'Energia'
$data[113,3] = '230'
$data[114,0] = 'Centrais de energia híbridas'
$data[114,1] = '234'
$data[114,2] = 'Energia'
$data[114,3] = '230'
$data[115,0] = 'Centrais de energia nuclear'
$data[115,1] = '235'
$data[115,2] = 'Energia'
$data[115,3] = '230'
$data[116,0] = 'Distribuição de energia'
$data[116,1] = '236'
$data[116,2] = 'Energia'
$data[116,3] = '230'
$data[117,0] = 'Distribuição de energia'
$data[117,1] = '236'
$data[117,2] = 'Energia'
$data[117,3] = '230'
$data[118,0] = 'Distribuição de energia'
$data[118,1] = '236'
$data[118,2] = 'Energia'
$data[118,3] = '230'
$data[119,0] = 'Distribuição de energia'
$data[119,1] = '236'
$data[119,2] = 'Energia'
$data[119,3] = '230'
$data[120,0] = 'Distribuição de energia'
$data[120,1] = '236'
$data[120,2] = 'Energia'
$data[120,3] = '230'
$data[121,0] = 'Distribuição de energia'
$data[121,1] = '236'
$data[121,2] = 'Energia'
$data[121,3] = '230'
$data[122,0] = 'Distribuição de energia'
$data[122,1] = '236'
$data[122,2] = 'Energia'
$data[122,3] = '230'
$data[123,0] = 'Serviços bancários e financeiros'
$data[123,1] = '240'
$data[123,2] = 'Serviços bancários e financeiros'
$data[123,3] = '240'
$data[124,0] = 'Serviços bancários e financeiros'
$data[124,1] = '240'
$data[124,2] = 'Serviços bancários e financeiros'
$data[124,3] = '240'
$data[125,0] = 'Serviços bancários e financeiros'
$data[125,1] = '240'
$data[125,2] = 'Serviços bancários e financeiros'
$data[125,3] = '240'
$data[126,0] = 'Serviços bancários e financeiros'
$data[126,1] = '240'
$data[126,2] = 'Serviços bancários e financeiros'
$data[126,3] = '240'
$data[127,0] = 'Serviços bancários e financeiros'
$data[127,1] = '240'
$data[127,2] = 'Serviços bancários e financeiros'
$data[127,3] = '240'
$data[128,0] = 'Serviços bancários e financeiros'
$data[128,1] = '240'
$data[128,2] = 'Serviços bancários e financeiros'
$data[128,3] = '240'
$data[129,0] = 'Negócios e outros serviços'
$data[129,1] = '250'
$data[129,2] = 'Negócios e outros serviços'
$data[129,3] = '250'
$data[130,0] = 'Negócios e outros serviços'
$data[130,1] = '250'
$data[130,2] = 'Negócios e outros serviços'
$data[130,3] = '250'
$data[131,0] = 'Negócios e outros serviços'
$data[131,1] = '250'
$data[131,2] = 'Negócios e outros serviços'
$data[131,3] = '250'
$data[132,0] = 'Negócios e outros serviços'
$data[132,1] = '250'
$data[132,2] = 'Negócios e outros serviços'
$data[132,3] = '250'
$data[133,0] = 'Agricultura'
$data[133,1] = '311'
$data[133,2] = 'Agricultura, Silvicultura, Pesca'
$data[133,3] = '310'
$data[134,0] = 'Agricultura'
$data[134,1] = '311'
$data[134,2] = 'Agricultura, Silvicultura, Pesca'
$data[134,3] = '310'
$data[135,0] = 'Agricultura'
$data[135,1] = '311'
$data[135,2] = 'Agricultura, Silvicultura, Pesca'
$data[135,3] = '310'
$data[136,0] = 'Agricultura'
$data[136,1] = '311'
$data[136,2] = 'Agricultura, Silvicultura, Pesca'
$data[136,3] = '310'
$data[137,0] = 'Agricultura'
$data[137,1] = '311'
$data[137,2] = 'Agricultura, Silvicultura, Pesca'
$data[137,3] = '310'
$data[138,0] = 'Agricultura'
$data[138,1] = '311'
$data[138,2] = 'Agricultura, Silvicultura, Pesca'
$data[138,3] = '310'
$data[139,0] = 'Agricultura'
$data[139,1] = '311'
$data[139,2] = 'Agricultura, Silvicultura, Pesca'
$data[139,3] = '310'
$data[140,0] = 'Agricultura'
$data[140,1] = '311'
$data[140,2] = 'Agricultura, Silvicultura, Pesca'
$data[140,3] = '310'
$data[141,0] = 'Agricultura'
$data[141,1] = '311'
$data[141,2] = 'Agricultura, Silvicultura, Pesca'
$data[141,3] = '310'
$data[142,0] = 'Agricultura'
$data[142,1] = '311'
$data[142,2] = 'Agricultura, Silvicultura, Pesca'
$data[142,3] = '310'
$data[143,0] = 'Agricultura'
$data[143,1] = '311'
$data[143,2] = 'Agricultura, Silvicultura, Pesca'
$data[143,3] = '310'
$data[144,0] = 'Agricultura'
$data[144,1] = '311'
$data[144,2] = 'Agricultura, Silvicultura, Pesca'
$data[144,3] = '310'
$data[145,0] = 'Agricultura'
$data[145,1] = '311'
$data[145,2] = 'Agricultura, Silvicultura, Pesca'
$data[145,3] = '310'
$data[146,0] = 'Agricultura'
$data[146,1] = '311'
$data[146,2] = 'Agricultura, Silvicultura, Pesca'
$data[146,3] = '310'
$data[147,0] = 'Agricultura'
$data[147,1] = '311'
$data[147,2] = 'Agricultura, Silvicultura, Pesca'
$data[147,3] = '310'
$data[148,0] = 'Agricultura'
$data[148,1] = '311'
$data[148,2] = 'Agricultura, Silvicultura, Pesca'
$data[148,3] = '310'
$data[149,0] = 'Agricultura'
$data[149,1] = '311'
$data[149,2] = 'Agricultura, Silvicultura, Pesca'
$data[149,3] = '310'
$data[150,0] = 'Agricultura'
$data[150,1] = '311'
$data[150,2] = 'Agricultura, Silvicultura, Pesca'
$data[150,3] = '310'
$data[151,0] = 'Silvicultura'
$data[151,1] = '312'
$data[151,2] = 'Agricultura, Silvicultura, Pesca'
$data[151,3] = '310'
$data[152,0] = 'Silvicultura'
$data[152,1] = '312'
$data[152,2] = 'Agricultura, Silvicultura, Pesca'
$data[152,3] = '310'
$data[153,0] = 'Silvicultura'
$data[153,1] = '312'
$data[153,2] = 'Agricultura, Silvicultura, Pesca'
$data[153,3] = '310'
$data[154,0] = 'Silvicultura'
$data[154,1] = '312'
$data[154,2] = 'Agricultura, Silvicultura, Pesca'
$data[154,3] = '310'
$data[155,0] = 'Silvicultura'
$data[155,1] = '312'
$data[155,2] = 'Agricultura, Silvicultura, Pesca'
$data[155,3] = '310'
$data[156,0] = 'Silvicultura'
$data[156,1] = '312'
$data[156,2] = 'Agricultura, Silvicultura, Pesca'
$data[156,3] = '310'
$data[157,0] = 'Pesca'
$data[157,1] = '313'
$data[157,2] = 'Agricultura, Silvicultura, Pesca'
$data[157,3] = '310'
$data[158,0] = 'Pesca'
$data[158,1] = '313'
$data[158,2] = 'Agricultura, Silvicultura, Pesca'
$data[158,3] = '310'
$data[159,0] = 'Pesca'
$data[159,1] = '313'
$data[159,2] = 'Agricultura, Silvicultura, Pesca'
$data[159,3] = '310'
$data[160,0] = 'Pesca'
$data[160,1] = '313'
$data[160,2] = 'Agricultura, Silvicultura, Pesca'
$data[160,3] = '310'
$data[161,0] = 'Pesca'
$data[161,1] = '313'
$data[161,2] = 'Agricultura, Silvicultura, Pesca'
$data[161,3] = '310'
$data[162,0] = 'Indústria'
$data[162,1] = '321'
$data[162,2] = 'Indústria, extractivas, construção'
$data[162,3] = '320'
$data[163,0] = 'Indústria'
$data[163,1] = '321'
$data[163,2] = 'Indústria, extractivas, construção'
$data[163,3] = '320'
$data[164,0] = 'Indústria'
$data[164,1] = '321'
$data[164,2] = 'Indústria, extractivas, construção'
$data[164,3] = '320'
$data[165,0] = 'Indústria'
$data[165,1] = '321'
$data[165,2] = 'Indústria, extractivas, construção'
$data[165,3] = '320'
$data[166,0] = 'Indústria'
$data[166,1] = '321'
$data[166,2] = 'Indústria, extractivas, construção'
$data[166,3] = '320'
$data[167,0] = 'Indústria'
$data[167,1] = '321'
$data[167,2] = 'Indústria, extractivas, construção'
$data[167,3] = '320'
$data[168,0] = 'Indústria'
$data[168,1] = '321'
$data[168,2] = 'Indústria, extractivas, construção'
$data[168,3] = '320'
$data[169,0] = 'Indústria'
$data[169,1] = '321'
$data[169,2] = 'Indústria, extractivas, construção'
$data[169,3] = '320'
$data[170,0] = 'Indústria'
$data[170,1] = '321'
$data[170,2] = 'Indústria, extractivas, construção'
$data[170,3] = '320'
$data[171,0] = 'Indústria'
$data[171,1] = '321'
$data[171,2] = 'Indústria, extractivas, construção'
$data[171,3] = '320'
$data[172,0] = 'Indústria'
$data[172,1] = '321'
$data[172,2] = 'Indústria, extractivas, construção'
$data[172,3] = '320'
$data[173,0] = 'Indústria'
$data[173,1] = '321'
$data[173,2] = 'Indústria, extractivas, construção'
$data[173,3] = '320'
$data[174,0] = 'Indústria'
$data[174,1] = '321'
$data[174,2] = 'Indústria, extractivas, construção'
$data[174,3] = '320'
$data[175,0] = 'Indústria'
$data[175,1] = '321'
$data[175,2] = 'Indústria, extractivas, construção'
$data[175,3] = '320'
$data[176,0] = 'Indústria'
$data[176,1] = '321'
$data[176,2] = 'Indústria, extractivas, construção'
$data[176,3] = '320'
$data[177,0] = 'Indústria'
$data[177,1] = '321'
$data[177,2] = 'Indústria, extractivas, construção'
$data[177,3] = '320'
$data[178,0] = 'Indústria'
$data[178,1] = '321'
$data[178,2] = 'Indústria, extractivas, construção'
$data[178,3] = '320'
$data[179,0] = 'Indústria'
$data[179,1] = '321'
$data[179,2] = 'Indústria, extractivas, construção'
$data[179,3] = '320'
$data[180,0] = 'Indústria'
$data[180,1] = '321'
$data[180,2] = 'Indústria, extractivas, construção'
$data[180,3] = '320'
$data[181,0] = 'Recursos minerais e mineração'
$data[181,1] = '322'
$data[181,2] = 'Indústria, extractivas, construção'
$data[181,3] = '320'
$data[182,0] = 'Recursos minerais e mineração'
$data[182,1] = '322'
$data[182,2] = 'Indústria, extractivas, construção'
$data[182,3] = '320'
$data[183,0] = 'Recursos minerais e mineração'
$data[183,1] = '322'
$data[183,2] = 'Indústria, extractivas, construção'
$data[183,3] = '320'
$data[184,0] = 'Recursos minerais e mineração'
$data[184,1] = '322'
$data[184,2] = 'Indústria, extractivas, construção'
$data[184,3] = '320'
$data[185,0] = 'Recursos minerais e mineração'
$data[185,1] = '322'
$data[185,2] = 'Indústria, extractivas, construção'
$data[185,3] = '320'
$data[186,0] = 'Recursos minerais e mineração'
$data[186,1] = '322'
$data[186,2] = 'Indústria, extractivas, construção'
$data[186,3] = '320'
$data[187,0] = 'Recursos minerais e mineração'
$data[187,1] = '322'
$data[187,2] = 'Indústria, extractivas, construção'
$data[187,3] = '320'
$data[188,0] = 'Recursos minerais e mineração'
$data[188,1] = '322'
$data[188,2] = 'Indústria, extractivas, construção'
$data[188,3] = '320'
$data[189,0] = 'Recursos minerais e mineração'
$data[189,1] = '322'
$data[189,2] = 'Indústria, extractivas, construção'
$data[189,3] = '320'
$data[190,0] = 'Recursos minerais e mineração'
$data[190,1] = '322'
$data[190,2] = 'Indústria, extractivas, construção'
$data[190,3] = '320'
$data[191,0] = 'Construção'
$data[191,1] = '323'
$data[191,2] = 'Indústria, extractivas, construção'
$data[191,3] = '320'
$data[192,0] = 'Políticas comerciais e regulamentos'
$data[192,1] = '331'
$data[192,2] = 'Políticas comerciais e regulamentos'
$data[192,3] = '331'
$data[193,0] = 'Políticas comerciais e regulamentos'
$data[193,1] = '331'
$data[193,2] = 'Políticas comerciais e regulamentos'
$data[193,3] = '331'
$data[194,0] = 'Políticas comerciais e regulamentos'
$data[194,1] = '331'
$data[194,2] = 'Políticas comerciais e regulamentos'
$data[194,3] = '331'
$data[195,0] = 'Políticas comerciais e regulamentos'
$data[195,1] = '331'
$data[195,2] = 'Políticas comerciais e regulamentos'
$data[195,3] = '331'
$data[196,0] = 'Políticas comerciais e regulamentos'
$data[196,1] = '331'
$data[196,2] = 'Políticas comerciais e regulamentos'
$data[196,3] = '331'
$data[197,0] = 'Políticas comerciais e regulamentos'
$data[197,1] = '331'
$data[197,2] = 'Políticas comerciais e regulamentos'
$data[197,3] = '331'
$data[198,0] = 'Turismo'
$data[198,1] = '332'
$data[198,2] = 'Turismo'
$data[198,3] = '332'
$data[199,0] = 'Proteção geral do ambiente'
$data[199,1] = '410'
$data[199,2] = 'Proteção geral do ambiente'
$data[199,3] = '410'
$data[200,0] = 'Proteção geral do ambiente'
$data[200,1] = '410'
$data[200,2] = 'Proteção geral do ambiente'
$data[200,3] = '410'
$data[201,0] = 'Proteção geral do ambiente'
$data[201,1] = '410'
$data[201,2] = 'Proteção geral do ambiente'
$data[201,3] = '410'
$data[202,0] = 'Proteção geral do ambiente'
$data[202,1] = '410'
$data[202,2] = 'Proteção geral do ambiente'
$data[202,3] = '410'
$data[203,0] = 'Proteção geral do ambiente'
$data[203,1] = '410'
$data[203,2] = 'Proteção geral do ambiente'
$data[203,3] = '410'
$data[204,0] = 'Proteção geral do ambiente'
$data[204,1] = '410'
$data[204,2] = 'Proteção geral do ambiente'
$data[204,3] = '410'
$data[205,0] = 'Outros, multissetoriais'
$data[205,1] = '430'
$data[205,2] = 'Outros, multissetoriais'
$data[205,3] = '430'
$data[206,0] = 'Outros, multissetoriais'
$data[206,1] = '430'
$data[206,2] = 'Outros, multissetoriais'
$data[206,3] = '430'
$data[207,0] = 'Outros, multissetoriais'
$data[207,1] = '430'
$data[207,2] = 'Outros, multissetoriais'
$data[207,3] = '430'
$data[208,0] = 'Outros, multissetoriais'
$data[208,1] = '430'
$data[208,2] = 'Outros, multissetoriais'
$data[208,3] = '430'
$data[209,0] = 'Outros, multissetoriais'
$data[209,1] = '430'
$data[209,2] = 'Outros, multissetoriais'
$data[209,3] = '430'
$data[210,0] = 'Outros, multissetoriais'
$data[210,1] = '430'
$data[210,2] = 'Outros, multissetoriais'
$data[210,3] = '430'
$data[211,0] = 'Outros, multissetoriais'
$data[211,1] = '430'
$data[211,2] = 'Outros, multissetoriais'
$data[211,3] = '430'
$data[212,0] = 'Outros, multissetoriais'
$data[212,1] = '430'
$data[212,2] = 'Outros, multissetoriais'
$data[212,3] = '430'
$data[213,0] = 'Outros, multissetoriais'
$data[213,1] = '430'
$data[213,2] = 'Outros, multissetoriais'
$data[213,3] = '430'
$data[214,0] = 'Outros, multissetoriais'
$data[214,1] = '430'
$data[214,2] = 'Outros, multissetoriais'
$data[214,3] = '430'
$data[215,0] = 'Apoio orçamental geral'
$data[215,1] = '510'
$data[215,2] = 'Apoio orçamental geral'
$data[215,3] = '510'
$data[216,0] = 'Ajuda alimentar ao desenvolvimento'
$data[216,1] = '520'
$data[216,2] = 'Ajuda alimentar ao desenvolvimento'
$data[216,3] = '520'
$data[217,0] = 'Outra ajuda por meio de mercadorias'
$data[217,1] = '530'
$data[217,2] = 'Outra ajuda por meio de mercadorias'
$data[217,3] = '530'
$data[218,0] = 'Outra ajuda por meio de mercadorias'
$data[218,1] = '530'
$data[218,2] = 'Outra ajuda por meio de mercadorias'
$data[218,3] = '530'
$data[219,0] = 'Ação relacionada com a dívida'
$data[219,1] = '600'
$data[219,2] = 'Ação relacionada com a dívida'
$data[219,3] = '600'
$data[220,0] = 'Ação relacionada com a dívida'
$data[220,1] = '600'
$data[220,2] = 'Ação relacionada com a dívida'
$data[220,3] = '600'
$data[221,0] = 'Ação relacionada com a dívida'
$data[221,1] = '600'
$data[221,2] = 'Ação relacionada com a dívida'
$data[221,3] = '600'
$data[222,0] = 'Ação relacionada com a dívida'
$data[222,1] = '600'
$data[222,2] = 'Ação relacionada com a dívida'
$data[222,3] = '600'
$data[223,0] = 'Ação relacionada com a dívida'
$data[223,1] = '600'
$data[223,2] = 'Ação relacionada com a dívida'
$data[223,3] = '600'
$data[224,0] = 'Ação relacionada com a dívida'
$data[224,1] = '600'
$data[224,2] = 'Ação relacionada com a dívida'
$data[224,3] = '600'
$data[225,0] = 'Ação relacionada com a dívida'
$data[225,1] = '600'
$data[225,2] = 'Ação relacionada com a dívida'
$data[225,3] = '600'
$data[226,0] = 'Resposta de emergência'
$data[226,1] = '720'
$data[226,2] = 'Resposta de emergência'
$data[226,3] = '720'
$data[227,0] = 'Resposta de emergência'
$data[227,1] = '720'
$data[227,2] = 'Resposta de emergência'
$data[227,3] = '720'
$data[228,0] = 'Resposta de emergência'
$data[228,1] = '720'
$data[228,2] = 'Resposta de emergência'
$data[228,3] = '720'
$data[229,0] = 'Ajuda à reconstrução e reabilitação'
$data[229,1] = '730'
$data[229,2] = 'Ajuda à reconstrução e reabilitação'
$data[229,3] = '730'
$data[230,0] = 'Prevenção e preparação contra catástrofes'
$data[230,1] = '740'
$data[230,2] = 'Prevenção e preparação contra catástrofes'
$data[230,3] = '740'
$data[231,0] = 'Custos administrativos dos doadores'
$data[231,1] = '910'
$data[231,2] = 'Custos administrativos dos doadores'
$data[231,3] = '910'
$data[232,0] = 'Refugiados em países doadores'
$data[232,1] = '930'
$data[232,2] = 'Refugiados em países doadores'
$data[232,3] = '930'
$data[233,0] = 'Não atribuído/não especificado'
$data[233,1] = '998'
$data[233,2] = 'Não atribuído/não especificado'
$data[233,3] = '998'
$data[234,0] = 'Não atribuído/não especificado'
$data[234,1] = '998'
$data[234,2] = 'Não atribuído/não especificado'
$data[234,3] = '998'
$ws.Range("D1:G235").Value = $data
